$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "legend" cells with shared strings (order matters so the
#     shared-string table is built with the same indices as the target:
#     6=(4; 1) 7=(5; 1) 8=(4; 2) 9=(5; 2)) ---

# 2x2 block V2:W3
$ws.Range("V2").Value = "(4; 1)"
$ws.Range("W2").Value = "(5; 1)"
$ws.Range("V3").Value = "(4; 2)"
$ws.Range("W3").Value = "(5; 2)"

# Single row Z2:AC2 (reuses the same strings)
$ws.Range("Z2").Value = "(4; 1)"
$ws.Range("AA2").Value = "(4; 2)"
$ws.Range("AB2").Value = "(5; 1)"
$ws.Range("AC2").Value = "(5; 2)"

# --- Highlight E2:F3 with the Accent 6 theme color (green). Apply the
#     theme color to one cell, then copy its format to the other three so
#     the engine only mints a single new style entry (matching the target
#     cellXfs layout) instead of one per cell. ---
$ws.Range("E2").Interior.ThemeColor = 10
$ws.Range("E2").Copy()
$ws.Range("F2:F3").PasteSpecial(-4122)
$ws.Range("E3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update the active selection to match the author's final cursor spot ---
$ws.Range("Z4").Select()
